$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 26: add I26/J26 timer checkpoint values ---
$ws.Cells.Item(26, 9).Value = 16182
$ws.Cells.Item(26, 10).Value = 16243

# --- New rows 56-70: World 5, level 5-3 / 5-G data ---

$ws.Cells.Item(56, 1).Value = "Enter 5-3"
$ws.Cells.Item(56, 2).Value = 18945
$ws.Cells.Item(56, 3).Value = 22185
$ws.Cells.Item(56, 4).Formula = "=IF(B56 >  0,C56-B56, 0)"

$ws.Cells.Item(57, 1).Value = "1st Move"
$ws.Cells.Item(57, 2).Value = 19177
$ws.Cells.Item(57, 3).Value = 22443
$ws.Cells.Item(57, 4).Formula = "=IF(B57 >  0,C57-B57, 0)"

$ws.Cells.Item(58, 1).Value = "Checkpoint 2676"
$ws.Cells.Item(58, 2).Value = 20120
$ws.Cells.Item(58, 3).Value = 23386
$ws.Cells.Item(58, 4).Formula = "=IF(B58 >  0,C58-B58, 0)"

$ws.Cells.Item(59, 1).Value = "Checkpoint 2869"
$ws.Cells.Item(59, 2).Value = 20187
$ws.Cells.Item(59, 3).Value = 23453
$ws.Cells.Item(59, 4).Formula = "=IF(B59 >  0,C59-B59, 0)"

$ws.Cells.Item(60, 1).Value = "Checkpoint 3080"
$ws.Cells.Item(60, 2).Value = 20258
$ws.Cells.Item(60, 3).Value = 23524
$ws.Cells.Item(60, 4).Formula = "=IF(B60 >  0,C60-B60, 0)"

$ws.Cells.Item(61, 1).Value = "Checkpoint 3355/3356"
$ws.Cells.Item(61, 2).Value = 20351
$ws.Cells.Item(61, 3).Value = 23618
$ws.Cells.Item(61, 4).Formula = "=IF(B61 >  0,C61-B61, 0)"

$ws.Cells.Item(62, 1).Value = "Checkpoint 3484/3485"
$ws.Cells.Item(62, 2).Value = 20394
$ws.Cells.Item(62, 3).Value = 23661
$ws.Cells.Item(62, 4).Formula = "=IF(B62 >  0,C62-B62, 0)"

$ws.Cells.Item(63, 1).Value = "Checkpoint 3692/3693"
$ws.Cells.Item(63, 2).Value = 20465
$ws.Cells.Item(63, 3).Value = 23733
$ws.Cells.Item(63, 4).Formula = "=IF(B63 >  0,C63-B63, 0)"

$ws.Cells.Item(64, 1).Value = "Get flag"
$ws.Cells.Item(64, 2).Value = 20658
$ws.Cells.Item(64, 3).Value = 23926
$ws.Cells.Item(64, 4).Formula = "=IF(B64 >  0,C64-B64, 0)"

$ws.Cells.Item(65, 1).Value = "Level end"
$ws.Cells.Item(65, 2).Value = 21172
$ws.Cells.Item(65, 3).Value = 24440
$ws.Cells.Item(65, 4).Formula = "=IF(B65 >  0,C65-B65, 0)"

$ws.Cells.Item(66, 1).Value = "Enter 5-G"
$ws.Cells.Item(66, 2).Value = 21661
$ws.Cells.Item(66, 3).Value = 25475
$ws.Cells.Item(66, 4).Formula = "=IF(B66 >  0,C66-B66, 0)"

$ws.Cells.Item(67, 1).Value = "1st move"
$ws.Cells.Item(67, 2).Value = 21890
$ws.Cells.Item(67, 3).Value = 25726
$ws.Cells.Item(67, 4).Formula = "=IF(B67 >  0,C67-B67, 0)"

$ws.Cells.Item(68, 1).Value = "Hit block"
$ws.Cells.Item(68, 2).Value = 21929
$ws.Cells.Item(68, 3).Value = 25765
$ws.Cells.Item(68, 4).Formula = "=IF(B68 >  0,C68-B68, 0)"

$ws.Cells.Item(69, 1).Value = "Checkpoint 10400"
$ws.Cells.Item(69, 2).Value = 21995
$ws.Cells.Item(69, 3).Value = 25832
$ws.Cells.Item(69, 4).Formula = "=IF(B69 >  0,C69-B69, 0)"

$ws.Cells.Item(70, 1).Value = "Get mini-mario"
$ws.Cells.Item(70, 3).Value = 25957

# --- View state: move selection to C71 (matches the author's final cursor
#     position after entering this new data) ---
$ws.Range("C71").Select()
